$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 712: mqtt_report_data_help -> split out "statistics" wording (device-info only now)
$ws.Cells.Item(712,2).Value = "Automatski šaljite podatke o uređaju putem MQTT svake 2 minute.<br><i>Za primjenu ove promjene potrebno je ponovno učitavanje.</i>"
$ws.Cells.Item(712,3).Value = "Automaticky zasílat informace o zařízení přes MQTT každé dvě minuty. <br><i>Restart aplikace je vyžadován pro aplikování tohoto nastavení.</i>"
$ws.Cells.Item(712,4).Value = "Automatically send device info through MQTT every 2 minutes.<br><i>Reload is required for applying this change.</i>"
$ws.Cells.Item(712,5).Value = "Envoyez automatiquement les informations sur l'appareil via MQTT toutes les 2 minutes.<br><i>Un rechargement est nécessaire pour appliquer cette modification.</i>"
$ws.Cells.Item(712,6).Value = "Senden Sie automatisch alle 2 Minuten Geräteinformationen über MQTT.<br><i>Neu laden ist erforderlich, um diese Änderung zu übernehmen.</i>"
$ws.Cells.Item(712,7).Value = "Invia automaticamente informazioni sul dispositivo tramite MQTT ogni 2 minuti.<br><i>Per applicare questa modifica è necessario ricaricare.</i>"
$ws.Cells.Item(712,8).Value = "Automatycznie wysyłaj informacje o urządzeniu przez MQTT co 2 minuty.<br><i>Do zastosowania tej zmiany wymagane jest ponowne załadowanie.</i>"
$ws.Cells.Item(712,9).Value = "Envie automaticamente informações do dispositivo por meio do MQTT a cada 2 minutos.<br><i>É necessário recarregar para aplicar essa alteração.</i>"
$ws.Cells.Item(712,10).Value = "Автоматически отправлять информацию об устройстве через MQTT каждые 2 минуты.<br><i>Для применения этого изменения требуется перезагрузка.</i>"
$ws.Cells.Item(712,11).Value = "Automaticky zasielať informácie o zariadení cez MQTT každé dve minúty. <br><i>Reštart aplikácie je vyžadovaný pre aplikovanie tohto nastavenia.</i>"
$ws.Cells.Item(712,12).Value = "Envíe automáticamente información del dispositivo a través de MQTT cada 2 minutos.<br><i>Es necesario volver a cargar para aplicar este cambio.</i>"

# Row 715
$ws.Cells.Item(715,1).Value = "mqtt_report_statistics"
$ws.Cells.Item(715,2).Value = "Prijavite statistiku putem MQTT-a"
$ws.Cells.Item(715,3).Value = "Zasílat statistiky přes MQTT"
$ws.Cells.Item(715,4).Value = "Report statistics through MQTT"
$ws.Cells.Item(715,5).Value = "Statistiques de rapport via MQTT"
$ws.Cells.Item(715,6).Value = "Melden Sie Statistiken über MQTT"
$ws.Cells.Item(715,7).Value = "Segnala le statistiche tramite MQTT"
$ws.Cells.Item(715,8).Value = "Raportuj statystyki przez MQTT"
$ws.Cells.Item(715,9).Value = "Relatório de estatísticas por meio do MQTT"
$ws.Cells.Item(715,10).Value = "Отчет о статистике через MQTT"
$ws.Cells.Item(715,11).Value = "Zasielať štatistiky cez MQTT"
$ws.Cells.Item(715,12).Value = "Informe de estadísticas a través de MQTT"

# Row 716
$ws.Cells.Item(716,1).Value = "mqtt_report_statistics_help"
$ws.Cells.Item(716,2).Value = "Izvještaj o statistici prikaza datoteka putem MQTT-a<br><i>Za primjenu ove promjene potrebno je ponovno učitavanje.</i>"
$ws.Cells.Item(716,3).Value = "Automaticky zasílat statistiky o zobrazených souborech přes MQTT každé dvě minuty.<br><i>Restart aplikace je vyžadován pro aplikování tohoto nastavení.</i>"
$ws.Cells.Item(716,4).Value = "Automatically report display statistics of files through MQTT every two minutes.<br><i>Reload is required for applying this change.</i>"
$ws.Cells.Item(716,5).Value = "Signaler les statistiques d'affichage des fichiers via MQTT<br><i>Un rechargement est nécessaire pour appliquer cette modification.</i>"
$ws.Cells.Item(716,6).Value = "Melden Sie Anzeigestatistiken von Dateien über MQTT<br><i>Neues Laden ist erforderlich, um diese Änderung zu übernehmen.</i>"
$ws.Cells.Item(716,7).Value = "Segnala le statistiche di visualizzazione dei file tramite MQTT<br><i>È necessario ricaricare per applicare questa modifica.</i>"
$ws.Cells.Item(716,8).Value = "Raportuj statystyki wyświetlania plików przez MQTT<br><i>Do zastosowania tej zmiany wymagane jest ponowne wczytanie.</i>"
$ws.Cells.Item(716,9).Value = "Relatório de estatísticas de exibição de arquivos por meio do MQTT<br><i>Recarregar é necessário para aplicar essa alteração.</i>"
$ws.Cells.Item(716,10).Value = "Отчет о статистике отображения файлов через MQTT<br><i>Для применения этого изменения требуется перезагрузка.</i>"
$ws.Cells.Item(716,11).Value = "Automaticky zasielať štatistiky o zobrazených súboroch cez MQTT každé dve minúty.<br><i>Reštart aplikácie je vyžadovaný pre aplikovanie tohto nastavenia.</i>"
$ws.Cells.Item(716,12).Value = "Informe de estadísticas de visualización de archivos a través de MQTT<br><i>Es necesario volver a cargar para aplicar este cambio.</i>"

# Row 717
$ws.Cells.Item(717,1).Value = "license_key_enter"
$ws.Cells.Item(717,2).Value = "Unesite licencni ključ"
$ws.Cells.Item(717,3).Value = "Zadat licenční klíč"
$ws.Cells.Item(717,4).Value = "Enter license key"
$ws.Cells.Item(717,5).Value = "Entrez la clé de licence"
$ws.Cells.Item(717,6).Value = "Gebe den Lizenzschlüssel ein"
$ws.Cells.Item(717,7).Value = "Inserisci la chiave di licenza"
$ws.Cells.Item(717,8).Value = "Wprowadź klucz licencyjny"
$ws.Cells.Item(717,9).Value = "Insira chave da licença"
$ws.Cells.Item(717,10).Value = "Введите лицензионный ключ"
$ws.Cells.Item(717,11).Value = "Zadať licenčný kľúč"
$ws.Cells.Item(717,12).Value = "Introduzca la clave de la licencia"

# Row 718
$ws.Cells.Item(718,1).Value = "free_trial"
$ws.Cells.Item(718,2).Value = "Free trial"
$ws.Cells.Item(718,3).Value = "Free trial"
$ws.Cells.Item(718,4).Value = "Free trial"
$ws.Cells.Item(718,5).Value = "Free trial"
$ws.Cells.Item(718,6).Value = "Free trial"
$ws.Cells.Item(718,7).Value = "Free trial"
$ws.Cells.Item(718,8).Value = "Free trial"
$ws.Cells.Item(718,9).Value = "Free trial"
$ws.Cells.Item(718,10).Value = "Free trial"
$ws.Cells.Item(718,11).Value = "Free trial"
$ws.Cells.Item(718,12).Value = "Free trial"

# Row 719
$ws.Cells.Item(719,1).Value = "licence_key_how_to_get"
$ws.Cells.Item(719,2).Value = "Kontaktirajte podršku kako biste dobili licencni ključ"
$ws.Cells.Item(719,3).Value = "Prosím kontaktujte podporu pro získání licenčního klíče"
$ws.Cells.Item(719,4).Value = "Please contact support to get the license key"
$ws.Cells.Item(719,5).Value = "Veuillez contacter le support pour obtenir la clé de licence"
$ws.Cells.Item(719,6).Value = "Bitte wenden Sie sich an den Support, um den Lizenzschlüssel zu erhalten"
$ws.Cells.Item(719,7).Value = "Si prega di contattare l'assistenza per ottenere la chiave di licenza"
$ws.Cells.Item(719,8).Value = "Skontaktuj się z pomocą techniczną, aby uzyskać klucz licencyjny"
$ws.Cells.Item(719,9).Value = "Para licenciar favor entrar em contato com o suporte"
$ws.Cells.Item(719,10).Value = "Обратитесь в службу поддержки, чтобы получить лицензионный ключ"
$ws.Cells.Item(719,11).Value = "Prosím kontaktujte podporu pre získanie licenčného kľúča"
$ws.Cells.Item(719,12).Value = "Póngase en contacto con el soporte para obtener la clave de licencia"

# Row 720
$ws.Cells.Item(720,1).Value = "license_key_device_id"
$ws.Cells.Item(720,2).Value = "Unesite licencni ključ ispod za sljedeći ID uređaja"
$ws.Cells.Item(720,3).Value = "Zadejte níže licenční klíč pro následující ID zařízení"
$ws.Cells.Item(720,4).Value = "Enter the license key below for the following device ID"
$ws.Cells.Item(720,5).Value = "Entrez la clé de licence ci-dessous pour l'ID d'appareil suivant"
$ws.Cells.Item(720,6).Value = "Geben Sie unten den Lizenzschlüssel für die folgende Geräte-ID ein"
$ws.Cells.Item(720,7).Value = "Immettere la chiave di licenza di seguito per il seguente ID dispositivo"
$ws.Cells.Item(720,8).Value = "Wprowadź klucz licencyjny poniżej dla następującego identyfikatora urządzenia"
$ws.Cells.Item(720,9).Value = "Digite a chave de licença abaixo para o seguinte ID do dispositivo"
$ws.Cells.Item(720,10).Value = "Введите лицензионный ключ ниже для следующего идентификатора устройства"
$ws.Cells.Item(720,11).Value = "Zadajte nižšie licenčný kľúč pre nasledujúce ID zariadenia"
$ws.Cells.Item(720,12).Value = "Ingrese la clave de licencia a continuación para el siguiente ID de dispositivo"

# Row 721
$ws.Cells.Item(721,1).Value = "license_key_valid"
$ws.Cells.Item(721,2).Value = "Ključ licence je važeći, aplikacija će se ponovno pokrenuti"
$ws.Cells.Item(721,3).Value = "Licenční klíč je správný, aplikace se restartuje"
$ws.Cells.Item(721,4).Value = "License key is valid, application will be restarted"
$ws.Cells.Item(721,5).Value = "La clé de licence est valide, l'application va être redémarrée"
$ws.Cells.Item(721,6).Value = "Der Lizenzschlüssel ist gültig, die Anwendung wird neu gestartet"
$ws.Cells.Item(721,7).Value = "La chiave di licenza è valida, l'applicazione verrà riavviata"
$ws.Cells.Item(721,8).Value = "Klucz licencyjny jest ważny, aplikacja zostanie ponownie uruchomiona"
$ws.Cells.Item(721,9).Value = "A chave de licença é válida, o aplicativo será reiniciado"
$ws.Cells.Item(721,10).Value = "Лицензионный ключ действителен, приложение будет перезапущено"
$ws.Cells.Item(721,11).Value = "Licenčný kľúč je správny, aplikácia sa reštartuje"
$ws.Cells.Item(721,12).Value = "La clave de licencia es válida, la aplicación se reiniciará"

# Row 722
$ws.Cells.Item(722,1).Value = "license_key_invalid"
$ws.Cells.Item(722,2).Value = "Ključ licence je nevažeći, provjerite je li ispravno unesen"
$ws.Cells.Item(722,3).Value = "Licenční klíč není správný, prosím zkontrolujte překlepy"
$ws.Cells.Item(722,4).Value = "License key is invalid, please check if it was entered correctly"
$ws.Cells.Item(722,5).Value = "La clé de licence n'est pas valide, veuillez vérifier si elle a été saisie correctement"
$ws.Cells.Item(722,6).Value = "Der Lizenzschlüssel ist ungültig, bitte überprüfen Sie, ob er korrekt eingegeben wurde"
$ws.Cells.Item(722,7).Value = "La chiave di licenza non è valida, controlla se è stata inserita correttamente"
$ws.Cells.Item(722,8).Value = "Klucz licencyjny jest nieprawidłowy, sprawdź, czy został wprowadzony poprawnie"
$ws.Cells.Item(722,9).Value = "A chave de licença é inválida, verifique se foi digitada corretamente"
$ws.Cells.Item(722,10).Value = "Лицензионный ключ недействителен, проверьте правильность ввода"
$ws.Cells.Item(722,11).Value = "Licenčný kľúč nie je správny, prosím skontrolujte preklepy"
$ws.Cells.Item(722,12).Value = "La clave de licencia no es válida, verifique si se ingresó correctamente"

$ws.Range("K721").Select()
